# Weekly fruit/vegetable price update: add one new week of Mandarina
# price rows at the top of the date-ordered block, pushing the older
# rows (which already existed) down by three rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the existing row 1029; this shifts
# every row from 1029 downward (through 1048) down to 1032-1051,
# preserving all of their existing data/formatting untouched.
$ws.Rows("1029:1031").Insert()

# --- New row 1029: Murcott / Especial, updated for the new week ---
$ws.Range("A1029").Value = 9
$ws.Range("B1029").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1029").Value = "Metropolitana"
$ws.Range("D1029").Value = 45239
$ws.Range("E1029").Value = 13
$ws.Range("F1029").Value = "Fruta"
$ws.Range("G1029").Value = 100102
$ws.Range("H1029").Value = "Cítricos"
$ws.Range("I1029").Value = 100102004
$ws.Range("J1029").Value = "Mandarina"
$ws.Range("K1029").Value = "Murcott"
$ws.Range("L1029").Value = "Especial"
$ws.Range("M1029").Value = 290
$ws.Range("N1029").Value = 8500
$ws.Range("O1029").Value = 8500
$ws.Range("P1029").Value = 8500
$ws.Range("Q1029").Value = "`$/bandeja 10 kilos"
$ws.Range("R1029").Value = "Provincia de Limarí"
$ws.Range("S1029").Value = 850
$ws.Range("T1029").Value = 10

# --- New row 1030: Murcott / Primera, updated for the new week ---
$ws.Range("A1030").Value = 9
$ws.Range("B1030").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1030").Value = "Metropolitana"
$ws.Range("D1030").Value = 45239
$ws.Range("E1030").Value = 13
$ws.Range("F1030").Value = "Fruta"
$ws.Range("G1030").Value = 100102
$ws.Range("H1030").Value = "Cítricos"
$ws.Range("I1030").Value = 100102004
$ws.Range("J1030").Value = "Mandarina"
$ws.Range("K1030").Value = "Murcott"
$ws.Range("L1030").Value = "Primera"
$ws.Range("M1030").Value = 300
$ws.Range("N1030").Value = 6500
$ws.Range("O1030").Value = 6500
$ws.Range("P1030").Value = 6500
$ws.Range("Q1030").Value = "`$/bandeja 10 kilos"
$ws.Range("R1030").Value = "Provincia de Limarí"
$ws.Range("S1030").Value = 650
$ws.Range("T1030").Value = 10

# --- New row 1031: Murcott / Segunda, updated for the new week ---
$ws.Range("A1031").Value = 9
$ws.Range("B1031").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1031").Value = "Metropolitana"
$ws.Range("D1031").Value = 45239
$ws.Range("E1031").Value = 13
$ws.Range("F1031").Value = "Fruta"
$ws.Range("G1031").Value = 100102
$ws.Range("H1031").Value = "Cítricos"
$ws.Range("I1031").Value = 100102004
$ws.Range("J1031").Value = "Mandarina"
$ws.Range("K1031").Value = "Murcott"
$ws.Range("L1031").Value = "Segunda"
$ws.Range("M1031").Value = 200
$ws.Range("N1031").Value = 4500
$ws.Range("O1031").Value = 4500
$ws.Range("P1031").Value = 4500
$ws.Range("Q1031").Value = "`$/bandeja 10 kilos"
$ws.Range("R1031").Value = "Provincia de Limarí"
$ws.Range("S1031").Value = 450
$ws.Range("T1031").Value = 10
